$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.07943533333333333
$ws.Range("H2").Value = 0.238306
$ws.Range("I2").Value = 0.8002511845635669
$ws.Range("J2").Value = 0.8002511845635669
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5706193333333334
$ws.Range("N2").Value = 1.711858
$ws.Range("O2").Value = 0.4188640502130462
$ws.Range("P2").Value = 0.4188640502130463
$ws.Range("Q2").Value = 0.04532733694977778
$ws.Range("R2").Value = 0.407946032548
$ws.Range("S2").Value = 0.3351964523540836
$ws.Range("T2").Value = 0.3351964523540837

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.07943533333333333
$ws.Range("H3").Value = 0.238306
$ws.Range("I3").Value = 0.8002511845635669
$ws.Range("J3").Value = 0.8002511845635669
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4846943333333333
$ws.Range("N3").Value = 1.454083
$ws.Range("O3").Value = 0.3557906641356566
$ws.Range("P3").Value = 0.3557906641356566
$ws.Range("Q3").Value = 0.0385018559331111
$ws.Range("R3").Value = 0.3465167033979999
$ws.Range("S3").Value = 0.2847219004312174
$ws.Range("T3").Value = 0.2847219004312174

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.07943533333333333
$ws.Range("H4").Value = 0.238306
$ws.Range("I4").Value = 0.8002511845635669
$ws.Range("J4").Value = 0.8002511845635669
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3069883333333334
$ws.Range("N4").Value = 0.920965
$ws.Range("O4").Value = 0.2253452856512971
$ws.Range("P4").Value = 0.2253452856512971
$ws.Range("Q4").Value = 0.02438572058777778
$ws.Range("R4").Value = 0.21947148529
$ws.Range("S4").Value = 0.1803328317782659
$ws.Range("T4").Value = 0.1803328317782659

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 0.01982766666666667
$ws.Range("H5").Value = 0.059483
$ws.Range("I5").Value = 0.1997488154364332
$ws.Range("J5").Value = 0.1997488154364332
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.5706193333333334
$ws.Range("N5").Value = 1.711858
$ws.Range("O5").Value = 0.4188640502130462
$ws.Range("P5").Value = 0.4188640502130463
$ws.Range("Q5").Value = 0.01131404993488889
$ws.Range("R5").Value = 0.101826449414
$ws.Range("S5").Value = 0.08366759785896266
$ws.Range("T5").Value = 0.08366759785896268

# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 0.01982766666666667
$ws.Range("H6").Value = 0.059483
$ws.Range("I6").Value = 0.1997488154364332
$ws.Range("J6").Value = 0.1997488154364332
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4846943333333333
$ws.Range("N6").Value = 1.454083
$ws.Range("O6").Value = 0.3557906641356566
$ws.Range("P6").Value = 0.3557906641356566
$ws.Range("Q6").Value = 0.009610357676555555
$ws.Range("R6").Value = 0.08649321908899998
$ws.Range("S6").Value = 0.07106876370443925
$ws.Range("T6").Value = 0.07106876370443926

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 0.01982766666666667
$ws.Range("H7").Value = 0.059483
$ws.Range("I7").Value = 0.1997488154364332
$ws.Range("J7").Value = 0.1997488154364332
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3069883333333334
$ws.Range("N7").Value = 0.920965
$ws.Range("O7").Value = 0.2253452856512971
$ws.Range("P7").Value = 0.2253452856512971
$ws.Range("Q7").Value = 0.00608686234388889
$ws.Range("R7").Value = 0.05478176109500001
$ws.Range("S7").Value = 0.04501245387303126
$ws.Range("T7").Value = 0.04501245387303126
